# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and "Correspond Handback
# DateTime" (col H) values for the row corresponding to
# c834d112-8d1f-43ed-b454-6a8a13c240f5.md (row 3) on both the "zh-cn" and
# "de-de" report sheets, reflecting the newly generated handback report
# timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-25 08:19:02"
$wsZhCn.Range("H3").Value = "2016-03-25 08:19:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-25 08:19:12"
$wsDeDe.Range("H3").Value = "2016-03-25 08:20:00"
